$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    37.78026849737784,
    48.72650743631324,
    53.38108332448111,
    69.18637462904245,
    78.7283180333721,
    89.86865947206783,
    104.9102670315499,
    114.4497776230204,
    127.9544756436364,
    137.8410833860163,
    144.3184152363771,
    159.9258231821239,
    169.288398864032,
    181.2185370051314,
    191.2606969940322,
    196.2424331534309,
    209.4036938034245,
    225.4779165271265,
    236.1211981883622,
    242.2162037202863,
    247.7348271188889,
    259.2852333583726,
    279.6313658886244,
    280.5491010327485,
    293.4673337497381,
    307.1278647175675,
    316.1046209500494,
    324.1343600576988,
    335.6226317676184,
    341.2372731593209,
    353.2226338130175,
    357.3131116996032,
    383.1703053500188,
    397.8856341758691,
    398.7776207906584,
    415.6449855088512,
    423.2452915393005,
    438.5885947038245,
    446.6896550694166,
    465.8083263651208,
    472.2497563846252,
    481.8761548945769,
    504.2771954897259,
    517.3957449828362,
    525.783650409262,
    539.5617452631528,
    540.7086053502101,
    557.0863253492348,
    564.1686228460293,
    572.0785584921365,
    581.89976530272,
    589.5159332108192,
    597.5251773286714,
    609.2398230932453,
    610.2495999744149,
    621.5724517240851,
    615.9212470505572,
    618.7569804071652,
    620.6545633920869,
    620.6463235810013,
    621.3295437090229
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

